$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.083.65'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.923.50'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.87'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4607'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3825'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9779'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.62'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.942.84'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.692'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.967'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07051'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.008'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '84.30'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009521'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.75'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.005'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.100.76'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.345'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.96'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.076'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.89'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '19.03'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.661'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '118.12'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.836'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09341'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.8556'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.120'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.243'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.026'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.160'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05684'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.170'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +17.73%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.006'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02046'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.11%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5515'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1756'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.315'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.000002828'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.197'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +6.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5192'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.06931'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.01%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.21'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '110.37'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.768'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.49%  '
